$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing values that changed (rows shifted up by "fixing" previous zero placeholders)
$ws.Range("C12").Value = 20
$ws.Range("C32").Value = 100
$ws.Range("C39").Value = 50
$ws.Range("C40").Value = 200

# Remove the last two rows (old rows 41 and 42), which are no longer present in the data
$ws.Rows("41:42").Delete()
